# dieu chinh hoa don co ma
# 1) Add two new paragraphs after the "<signNameSubTitle3>" paragraph
#    (last cell of the signer table): an empty TableParagraph, followed
#    by a TableParagraph holding the "_GoBack" bookmark and the
#    "<digitalSignature>" placeholder text.
# 2) Remove the "_GoBack" bookmark that used to live in the header
#    (first-page header) leading paragraph, leaving the paragraph itself
#    (and its spacing) intact.

$d = $word.ActiveDocument

# --- Part 1: body - insert the two new paragraphs -----------------------
$bodyRng = $d.Content
$found = $bodyRng.Find.Execute("<signNameSubTitle3>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find <signNameSubTitle3> placeholder"
}
$targetPara = $bodyRng.Paragraphs(1)

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="TableParagraph"/><w:jc w:val="center"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="TableParagraph"/><w:jc w:val="center"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&lt;digitalSignature&gt;</w:t></w:r></w:p>'

$insertPoint = $d.Range($targetPara.Range.End, $targetPara.Range.End)
$insertPoint.InsertXML($newParagraphsXml)

# --- Part 2: header - drop the old "_GoBack" bookmark -------------------
$sec = $d.Sections.Item(1)
$header = $sec.Headers.Item(2)   # wdHeaderFooterFirstPage -> header3.xml
$headerRng = $header.Range

# Re-insert the leading paragraph's properties as a fresh paragraph
# (without the bookmark) right before the original one, then delete the
# original paragraph mark so only the bookmark-free copy remains.
$headerRng.SetRange($headerRng.Start, $headerRng.Start)
$cleanParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="005A002D" w:rsidRDefault="005A002D" w:rsidP="00687F97"><w:pPr><w:spacing w:after="240"/></w:pPr></w:p>'
$headerRng.InsertXML($cleanParaXml)

$oldMark = $header.Range
$oldMark.SetRange(1, 2)
$oldMark.Delete()

Write-Host "edit complete"
